$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF") with the same style as the other headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data for rows 2..86: column I = "I0", column J = "IF"
$data = @(
        @(8,8),
        @(8,8),
        @(8,8),
        @(8,8),
        @(6,7),
        @(7,7),
        @(7,7),
        @(8,9),
        @(8,8),
        @(9,9),
        @(8,8),
        @(7,7),
        @(7,7),
        @(7,7),
        @(9,9),
        @(8,8),
        @(7,7),
        @(8,8),
        @(8,8),
        @(8,8),
        @(9,9),
        @(8,8),
        @(8,8),
        @(8,8),
        @(8,8),
        @(8,8),
        @(8,8),
        @(8,8),
        @(10,10),
        @(10,10),
        @(9,10),
        @(7,8),
        @(8,8),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(8,9),
        @(8,9),
        @(8,8),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(8,9),
        @(8,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(10,10),
        @(10,10),
        @(9,9),
        @(9,9),
        @(8,9),
        @(8,8),
        @(8,8),
        @(10,10),
        @(9,9),
        @(10,10),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(8,9),
        @(8,8),
        @(8,8),
        @(10,10),
        @(9,9),
        @(5,6),
        @(3,3),
        @(8,8),
        @(6,6),
        @(6,6),
        @(5,5),
        @(7,7),
        @(3,3)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $ws.Cells.Item($rowNum, 9).Value = $data[$r][0]
    $ws.Cells.Item($rowNum, 10).Value = $data[$r][1]
}
